$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header for the added date column (X) - "6-jul"
$ws.Range("X1").Value = "6-jul"

# New data values for column X (added next to existing W "5-jul" column)
$ws.Range("X2").Value = 0
$ws.Range("X3").Value = 13.696895105679078
$ws.Range("X4").Value = 11.878435897576063
$ws.Range("X5").Value = 25.839536848412195
$ws.Range("X6").Value = 0
$ws.Range("X7").Value = 26.660891473203652
$ws.Range("X8").Value = 13.049338848802693
$ws.Range("X9").Value = 21.831898765320943
$ws.Range("X10").Value = 25.005500657833284
$ws.Range("X11").Value = 12.024928644576905
$ws.Range("X12").Value = 0
$ws.Range("X13").Value = 14.464628206568607
$ws.Range("X14").Value = 0
$ws.Range("X15").Value = 0
$ws.Range("X16").Value = 13.885079242215715
$ws.Range("X17").Value = 0
$ws.Range("X18").Value = 0

# Apply the same number formatting (2 decimals) and center alignment used
# for the other date columns to the existing "5-jul" column (W), matching
# the newly introduced style.
$rngW = $ws.Range("W2:W18")
$rngW.HorizontalAlignment = -4108
$rngW.NumberFormat = "0.00"

# Update the active selection to the newly added column, as recorded by
# the workbook after the edit.
[void]$ws.Range("X2:X18").Select()
